# Primeros cambios despues de usar la apk
# Reworks "Hoja 1": new column order (name, tipo, peso, price, stock,
# fechaActualizacion, id) and 4 data rows (previously 1).
#
# Values that look like plain numbers or ISO dates are entered with a
# leading apostrophe so Excel keeps them as text (matching the sheet's
# "number stored as text" convention) instead of silently converting them
# to Number/Date. The "id" column is a genuine number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "name"
$ws.Cells.Item(1,2).Value = "tipo"
$ws.Cells.Item(1,3).Value = "peso"
$ws.Cells.Item(1,4).Value = "price"
$ws.Cells.Item(1,5).Value = "stock"
$ws.Cells.Item(1,6).Value = "fechaActualizacion"
$ws.Cells.Item(1,7).Value = "id"

# Row 2 - DOGUI
$ws.Cells.Item(2,1).Value = "DOGUI"
$ws.Cells.Item(2,2).Value = "PERRO"
$ws.Cells.Item(2,3).Value = "'23"
$ws.Cells.Item(2,4).Value = "'23"
$ws.Cells.Item(2,5).Value = "'23"
$ws.Cells.Item(2,6).Value = "'2023-09-07"
$ws.Cells.Item(2,7).Value = 1

# Row 3 - NUTRIBON
$ws.Cells.Item(3,1).Value = "NUTRIBON"
$ws.Cells.Item(3,2).Value = "PERRO"
$ws.Cells.Item(3,3).Value = "'20"
$ws.Cells.Item(3,4).Value = "'31"
$ws.Cells.Item(3,5).Value = "'31"
$ws.Cells.Item(3,6).Value = "'2023-09-07"
$ws.Cells.Item(3,7).Value = 2

# Row 4 - A
$ws.Cells.Item(4,1).Value = "A"
$ws.Cells.Item(4,2).Value = "PERRO"
$ws.Cells.Item(4,3).Value = "'123"
$ws.Cells.Item(4,4).Value = "'123"
$ws.Cells.Item(4,5).Value = "'123"
$ws.Cells.Item(4,6).Value = "'2023-09-07"
$ws.Cells.Item(4,7).Value = 3

# Row 5 - B
$ws.Cells.Item(5,1).Value = "B"
$ws.Cells.Item(5,2).Value = "PERRO"
$ws.Cells.Item(5,3).Value = "'123"
$ws.Cells.Item(5,4).Value = "'123"
$ws.Cells.Item(5,5).Value = "'123"
$ws.Cells.Item(5,6).Value = "'2023-09-07"
$ws.Cells.Item(5,7).Value = 4
